# Append Serie A Round 22 results (rows 210-219) to the match log sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 210: Torino vs Cagliari
# Copy formatting (bold/bordered style on col A) from the last existing row first.
$ws.Range("A209:O209").Copy($ws.Range("A210:O210"))
$ws.Range("A210").Value = 208
$ws.Range("B210").Value = "Torino"
$ws.Range("C210").Value = "Cagliari"
$ws.Range("D210").Value = 2
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 1.56
$ws.Range("G210").Value = 0.34
$ws.Range("H210").Value = 1.76
$ws.Range("I210").Value = 0.33
$ws.Range("J210").Value = 0
$ws.Range("K210").Value = 0
$ws.Range("L210").Value = 0.2
$ws.Range("M210").Value = 0.01
$ws.Range("N210").Value = 0.22
$ws.Range("O210").Value = 2

$ws.Range("A210:O210").Copy($ws.Range("A211:O211"))
$ws.Range("A211").Value = 209
$ws.Range("B211").Value = "Como"
$ws.Range("C211").Value = "Atalanta"
$ws.Range("D211").Value = 1
$ws.Range("E211").Value = 2
$ws.Range("F211").Value = 0.77
$ws.Range("G211").Value = 0.96
$ws.Range("H211").Value = 1.11
$ws.Range("I211").Value = 0.9
$ws.Range("J211").Value = 0
$ws.Range("K211").Value = 0
$ws.Range("L211").Value = 0.34
$ws.Range("M211").Value = 0.06
$ws.Range("N211").Value = 0.41
$ws.Range("O211").Value = 3

$ws.Range("A211:O211").Copy($ws.Range("A212:O212"))
$ws.Range("A212").Value = 210
$ws.Range("B212").Value = "Empoli"
$ws.Range("C212").Value = "Bologna"
$ws.Range("D212").Value = 1
$ws.Range("E212").Value = 1
$ws.Range("F212").Value = 0.57
$ws.Range("G212").Value = 0.67
$ws.Range("H212").Value = 0.92
$ws.Range("I212").Value = 1
$ws.Range("J212").Value = 1
$ws.Range("K212").Value = 0
$ws.Range("L212").Value = 0.35
$ws.Range("M212").Value = 0.33
$ws.Range("N212").Value = 0.68
$ws.Range("O212").Value = 1

$ws.Range("A212:O212").Copy($ws.Range("A213:O213"))
$ws.Range("A213").Value = 211
$ws.Range("B213").Value = "Napoli"
$ws.Range("C213").Value = "Juventus"
$ws.Range("D213").Value = 2
$ws.Range("E213").Value = 1
$ws.Range("F213").Value = 2.37
$ws.Range("G213").Value = 0.8100000000000001
$ws.Range("H213").Value = 2.88
$ws.Range("I213").Value = 0.55
$ws.Range("J213").Value = 2
$ws.Range("K213").Value = 0
$ws.Range("L213").Value = 0.51
$ws.Range("M213").Value = 0.26
$ws.Range("N213").Value = 0.77
$ws.Range("O213").Value = 1

$ws.Range("A213:O213").Copy($ws.Range("A214:O214"))
$ws.Range("A214").Value = 212
$ws.Range("B214").Value = "Lazio"
$ws.Range("C214").Value = "Fiorentina"
$ws.Range("D214").Value = 1
$ws.Range("E214").Value = 2
$ws.Range("F214").Value = 1.41
$ws.Range("G214").Value = 0.73
$ws.Range("H214").Value = 1.63
$ws.Range("I214").Value = 0.88
$ws.Range("J214").Value = 0
$ws.Range("K214").Value = 0
$ws.Range("L214").Value = 0.22
$ws.Range("M214").Value = 0.15
$ws.Range("N214").Value = 0.38
$ws.Range("O214").Value = 3

$ws.Range("A214:O214").Copy($ws.Range("A215:O215"))
$ws.Range("A215").Value = 213
$ws.Range("B215").Value = "Lecce"
$ws.Range("C215").Value = "Inter"
$ws.Range("D215").Value = 0
$ws.Range("E215").Value = 4
$ws.Range("F215").Value = 0.88
$ws.Range("G215").Value = 2.12
$ws.Range("H215").Value = 1.05
$ws.Range("I215").Value = 2.24
$ws.Range("J215").Value = 0
$ws.Range("K215").Value = 1
$ws.Range("L215").Value = 0.17
$ws.Range("M215").Value = 0.12
$ws.Range("N215").Value = 0.28
$ws.Range("O215").Value = 3

$ws.Range("A215:O215").Copy($ws.Range("A216:O216"))
$ws.Range("A216").Value = 214
$ws.Range("B216").Value = "Milan"
$ws.Range("C216").Value = "Parma"
$ws.Range("D216").Value = 3
$ws.Range("E216").Value = 2
$ws.Range("F216").Value = 2.32
$ws.Range("G216").Value = 2.33
$ws.Range("H216").Value = 2.55
$ws.Range("I216").Value = 1.43
$ws.Range("J216").Value = 2
$ws.Range("K216").Value = 0
$ws.Range("L216").Value = 0.23
$ws.Range("M216").Value = 0.9
$ws.Range("N216").Value = 1.12
$ws.Range("O216").Value = 3

$ws.Range("A216:O216").Copy($ws.Range("A217:O217"))
$ws.Range("A217").Value = 215
$ws.Range("B217").Value = "Udinese"
$ws.Range("C217").Value = "Roma"
$ws.Range("D217").Value = 1
$ws.Range("E217").Value = 2
$ws.Range("F217").Value = 0.8
$ws.Range("G217").Value = 2.85
$ws.Range("H217").Value = 0.48
$ws.Range("I217").Value = 3.46
$ws.Range("J217").Value = 0
$ws.Range("K217").Value = 2
$ws.Range("L217").Value = 0.32
$ws.Range("M217").Value = 0.61
$ws.Range("N217").Value = 0.92
$ws.Range("O217").Value = 1

$ws.Range("A217:O217").Copy($ws.Range("A218:O218"))
$ws.Range("A218").Value = 216
$ws.Range("B218").Value = "Genoa"
$ws.Range("C218").Value = "Monza"
$ws.Range("D218").Value = 2
$ws.Range("E218").Value = 0
$ws.Range("F218").Value = 2.66
$ws.Range("G218").Value = 0.43
$ws.Range("H218").Value = 2.41
$ws.Range("I218").Value = 0.43
$ws.Range("J218").Value = 1
$ws.Range("K218").Value = 0
$ws.Range("L218").Value = 0.25
$ws.Range("M218").Value = 0
$ws.Range("N218").Value = 0.25
$ws.Range("O218").Value = 1

$ws.Range("A218:O218").Copy($ws.Range("A219:O219"))
$ws.Range("A219").Value = 217
$ws.Range("B219").Value = "Venezia"
$ws.Range("C219").Value = "Hellas Verona"
$ws.Range("D219").Value = 1
$ws.Range("E219").Value = 1
$ws.Range("F219").Value = 0.58
$ws.Range("G219").Value = 1.76
$ws.Range("H219").Value = 0.76
$ws.Range("I219").Value = 1.99
$ws.Range("J219").Value = 0
$ws.Range("K219").Value = 0
$ws.Range("L219").Value = 0.18
$ws.Range("M219").Value = 0.23
$ws.Range("N219").Value = 0.41
$ws.Range("O219").Value = 2

